# udi's game - updated graphics
#
# Appends 7 new rows (115-121) to the "Users" sheet, repeating the same
# data pattern already present in rows 113-114:
#   A: moses   B: bro   C: "1234" (text)   D: m@g.c   E: Male   F: 0 (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$firstRow = 115
$lastRow = 121

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("A$r").Value = "moses"
    $ws.Range("B$r").Value = "bro"
    $ws.Range("D$r").Value = "m@g.c"
    $ws.Range("E$r").Value = "Male"
    $ws.Range("F$r").Value = 0
}

# Column C holds the digit string "1234" stored as TEXT (matching the
# existing rows above), not the number 1234. Force text entry via a
# temporary "@" (text) number format, write the value, then restore the
# "Normal" style so the cells end up with the same default style as the
# rest of the sheet.
$cRange = $ws.Range("C$firstRow`:C$lastRow")
$cRange.NumberFormat = "@"
$cRange.Value = "1234"
$cRange.Style = "Normal"
